$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new data row for the November 2024 transect
$ws.Range("A34").Value = 2024
$ws.Range("B34").Value = 11
$ws.Range("C34").Value = 15
$ws.Range("D34").Value = 1
$ws.Range("E34").Value = 1
$ws.Range("F34").Value = 1
$ws.Range("G34").Value = 0
$ws.Range("H34").Value = "presence_abscence"

# Update selection to reflect where the user ended up after entering data
$ws.Range("C35").Select()
